$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.202.50"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.313.78"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'256.71"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'629.18"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  +21.61%  "
$ws.Range("D8").Value = "'0.410"
$ws.Range("E8").Value = "  +6.64%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.984"
$ws.Range("E10").Value = "  +22.65%  "
$ws.Range("D11").Value = "3.313.23"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "'41.11"
$ws.Range("E13").Value = "  +14.87%  "
$ws.Range("D14").Value = "98.888.74"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'0.0000252"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "3.926.21"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "'5.50"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "3.305.60"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'3.50"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").Value = "'15.79"
$ws.Range("E20").Value = "  +7.23%  "
$ws.Range("D21").Value = "'6.43"
$ws.Range("E21").Value = "  +9.53%  "
$ws.Range("D22").Value = "'488.96"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'9.50"
$ws.Range("E23").Value = "  +4.43%  "
$ws.Range("D24").Value = "'0.0000205"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "'5.76"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  +36.72%  "
$ws.Range("D27").Value = "'89.16"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").Value = "'12.21"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").Value = "3.470.56"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "'0.152"
$ws.Range("E30").Value = "  +22.89%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'0.191"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'10.70"
$ws.Range("E33").Value = "  +16.82%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'28.10"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").Value = "'0.487"
$ws.Range("E36").Value = "  +8.91%  "
$ws.Range("D37").Value = "'0.153"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").Value = "'7.41"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").Value = "'1.97"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").Value = "'499.66"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "'3.85"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("D43").Value = "'1.26"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "'0.791"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'3.22"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'1.97"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").Value = "'159.24"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "'4.86"
$ws.Range("E49").Value = "  +7.77%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.856"
$ws.Range("E50").Value = "  +8.41%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'7.36"
$ws.Range("E51").Value = "  +15.81%  "
